$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.100.91"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "2.021.81"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.00"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.609"
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.13"
$ws.Range("E8").Value = "  -3.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.381"
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0787"
$ws.Range("E10").Value = "  +2.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("E11").Value = "  -3.06%  "
$ws.Range("D12").Value = "2.320.59"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.32"
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.52"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.744"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.15"
$ws.Range("D17").Value = "2.023.89"
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D18").Value = "37.039.27"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.17"
$ws.Range("E19").Value = "  +3.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.95"
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("D21").Value = "0.0₃0825"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.81"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("E24").Value = "  +3.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("E25").Value = "  -5.29%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.26"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.48"
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.126"
$ws.Range("E28").Value = "  -3.08%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.36"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.79"
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("E31").Value = "  -3.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.53"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0617"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.44"
$ws.Range("E34").Value = "  -2.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.36"
$ws.Range("E35").Value = "  -4.47%  "
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.16"
$ws.Range("E38").Value = "  -3.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.39"
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("E40").Value = "  -3.38%  "
$ws.Range("D41").Value = "1.478.86"
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.77"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "95.50"
$ws.Range("E43").Value = "  -2.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0925"
$ws.Range("E44").Value = "  -2.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.79"
$ws.Range("E45").Value = "  -3.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.14"
$ws.Range("E46").Value = "  -3.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.31"
$ws.Range("E47").Value = "  +1.61%  "
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").Value = "2.208.46"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.51"
$ws.Range("E51").Value = "  -1.49%  "
